$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Present-Storage")
$ws.Cells.Item(2, 2).Value = 132.04
$ws.Cells.Item(3, 2).Value = 563.3
$ws.Cells.Item(4, 2).Value = 58.67
$ws.Cells.Item(6, 2).Value = 139.11
$ws.Cells.Item(7, 2).Value = 595.17
$ws.Cells.Item(8, 2).Value = 61.52
$ws.Cells.Item(10, 2).Value = 107.99
$ws.Cells.Item(11, 2).Value = 444.46
$ws.Cells.Item(12, 2).Value = 50.75
$ws.Cells.Item(14, 2).Value = 128.01
$ws.Cells.Item(15, 2).Value = 541.8
$ws.Cells.Item(16, 2).Value = 57.61
$ws.Columns.Item(2).ColumnWidth = 8.8

$ws = $wb.Worksheets.Item("2030-Storage")
$ws.Cells.Item(2, 2).Value = 5.88
$ws.Cells.Item(3, 2).Value = 6.04
$ws.Cells.Item(4, 2).Value = 6.37
$ws.Cells.Item(6, 2).Value = 5.7
$ws.Cells.Item(7, 2).Value = 5.87
$ws.Cells.Item(8, 2).Value = 6.21
$ws.Cells.Item(10, 2).Value = 9.56
$ws.Cells.Item(11, 2).Value = 9.69
$ws.Cells.Item(12, 2).Value = 9.94
$ws.Cells.Item(14, 2).Value = 6.96
$ws.Cells.Item(15, 2).Value = 7.12
$ws.Cells.Item(16, 2).Value = 7.43

$ws = $wb.Worksheets.Item("2050-Storage")
$ws.Cells.Item(2, 2).Value = 146.2
$ws.Cells.Item(3, 2).Value = 73.36
$ws.Cells.Item(4, 2).Value = 92.23
$ws.Cells.Item(6, 2).Value = 154.08
$ws.Cells.Item(7, 2).Value = 77.05
$ws.Cells.Item(8, 2).Value = 97.01
$ws.Cells.Item(10, 2).Value = 119.03
$ws.Cells.Item(11, 2).Value = 62.2
$ws.Cells.Item(12, 2).Value = 76.93
$ws.Cells.Item(14, 2).Value = 141.59
$ws.Cells.Item(15, 2).Value = 71.7
$ws.Cells.Item(16, 2).Value = 89.81
$ws.Columns.Item(2).ColumnWidth = 8.8
